# Apply attendance counts to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value (cells previously held 0)
$updates = @{
    "G3"  = 1
    "H3"  = 1

    "D4"  = 2
    "E4"  = 1
    "F4"  = 1

    "D5"  = 1
    "E5"  = 1

    "D6"  = 1
    "E6"  = 1

    "H7"  = 1

    "H8"  = 1

    "D9"  = 1
    "E9"  = 1

    "G10" = 1
    "H10" = 1

    "H11" = 1

    "H12" = 1

    "H13" = 1

    "H14" = 1

    "D15" = 1
    "E15" = 1

    "H16" = 1

    "H17" = 1

    "H18" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$wb.Save()
